$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price/volume cells so numeric-looking strings stay as text
$textCells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'E13', 'E14', 'D15', 'E15', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'E25', 'D26', 'E26', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E46', 'E47', 'D49', 'E49', 'D50', 'E50', 'E51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '331.87'
$ws.Range('E2').Value = '1.00%'
$ws.Range('D3').Value = '45.87'
$ws.Range('E3').Value = '4.29%'
$ws.Range('D4').Value = '5.633'
$ws.Range('E4').Value = '2.18%'
$ws.Range('D5').Value = '0.08361'
$ws.Range('E5').Value = '4.24%'
$ws.Range('D6').Value = '2.041'
$ws.Range('E6').Value = '2.98%'
$ws.Range('D7').Value = '0.9750'
$ws.Range('E7').Value = '2.65%'
$ws.Range('E8').Value = '-0.87%'
$ws.Range('D9').Value = '0.1159'
$ws.Range('E9').Value = '1.83%'
$ws.Range('D10').Value = '0.1917'
$ws.Range('E10').Value = '1.64%'
$ws.Range('D11').Value = '10.37'
$ws.Range('E11').Value = '-2.79%'
$ws.Range('D12').Value = '0.09987'
$ws.Range('E12').Value = '0.26%'
$ws.Range('E13').Value = '-1.06%'
$ws.Range('E14').Value = '-0.39%'
$ws.Range('D15').Value = '0.001293'
$ws.Range('E15').Value = '1.39%'
$ws.Range('E16').Value = '1.29%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.377'
$ws.Range('E17').Value = '0.40%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '4.451'
$ws.Range('E18').Value = '1.87%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3364'
$ws.Range('E19').Value = '-3.16%'
$ws.Range('B20').Value = 'ProBitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D20').Value = '0.1392'
$ws.Range('E20').Value = '-1.88%'
$ws.Range('B21').Value = 'ZBToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D21').Value = '0.2652'
$ws.Range('E21').Value = '4.14%'
$ws.Range('B22').Value = 'CoinExToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D22').Value = '0.04190'
$ws.Range('E22').Value = '2.88%'
$ws.Range('B23').Value = 'BitKan'
$ws.Range('C23').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D23').Value = '0.001312'
$ws.Range('E23').Value = '3.77%'
$ws.Range('B24').Value = 'HotbitToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D24').Value = '0.004598'
$ws.Range('E24').Value = '6.08%'
$ws.Range('E25').Value = '8.49%'
$ws.Range('D26').Value = '0.0003746'
$ws.Range('E26').Value = '0.03%'
$ws.Range('D38').Value = '0.02767'
$ws.Range('E38').Value = '6.57%'
$ws.Range('D39').Value = '0.05812'
$ws.Range('E39').Value = '2.39%'
$ws.Range('D40').Value = '0.007750'
$ws.Range('E40').Value = '2.60%'
$ws.Range('D41').Value = '0.1438'
$ws.Range('E41').Value = '2.78%'
$ws.Range('D42').Value = '0.007271'
$ws.Range('E42').Value = '-2.86%'
$ws.Range('D43').Value = '0.002117'
$ws.Range('E43').Value = '5.05%'
$ws.Range('D44').Value = '0.008082'
$ws.Range('E44').Value = '-6.32%'
$ws.Range('D45').Value = '0.3405'
$ws.Range('E46').Value = '2.76%'
$ws.Range('E47').Value = '0.14%'
$ws.Range('D49').Value = '0.003506'
$ws.Range('E49').Value = '-7.00%'
$ws.Range('D50').Value = '0.003504'
$ws.Range('E50').Value = '-0.71%'
$ws.Range('E51').Value = '0.14%'
